$d = $word.ActiveDocument

# Short Summary (paragraph 3)
$d.Paragraphs(3).Range.Text = 'The article is about the rapid adoption of ChatGPT AI, which has a million unique users within five days of its release, causing discussions about the future of work in an AI-powered world.'

# Summary (paragraph 5)
$d.Paragraphs(5).Range.Text = '1. ChatGPT, an artificial intelligence large language model, has seen unprecedented adoption within just five days of its release.' + [char]11 + '2. Despite the rapid acceptance, the technology is still developing and lacks the ability to reason like humans.' + [char]11 + '3. The human mind uses a small amount of data to create thoughts, language, and the ability to reason, while AI models analyze massive datasets to produce content based on trends.' + [char]11 + '4. There are ethical concerns raised about the quick embrace of AI technologies, including their potential to reinforce ideologies and lack of balance between creativity and constraint.'

# Question 1 (paragraph 7)
$d.Paragraphs(7).Range.Text = 'The metaphors used in the media to frame the public discussion about ChatGPT include:' + [char]11 + [char]11 + '1. "Turning the world upside down" - This metaphor suggests that ChatGPT has caused a significant and transformative shift in society.' + [char]11 + '2. "AI age" or "start of the AI age" - This metaphor implies that ChatGPT represents the beginning of an era where artificial intelligence plays a dominant role in our lives.' + [char]11 + '3. "Predictable human behavior" - The media uses this metaphor to explain how platforms like ChatGPT analyze large data sets and make predictions about what users are thinking or looking for, suggesting that human behavior is becoming more predictable due to technology usage.' + [char]11 + '4. "The iPhone to TikTok" - This metaphor refers to the rapid adoption rates of new technologies, comparing ChatGPT to well-known examples like the iPhone and TikTok.' + [char]11 + '5. "Ethical issues" or "ethics debates" - This metaphor highlights concerns about the potential misuse and consequences of AI tools, such as reinforcing ideologies and worldviews, lacking creativity, and balancing creativity with constraint.' + [char]11 + '6. "Collaboration instead of a complete takeover or outsourcing" - This metaphor suggests that while AI tools can assist in various tasks, they cannot replace the human mind entirely and should be used collaboratively rather than as a replacement for human intelligence.'

# Question 2 (paragraph 9)
$d.Paragraphs(9).Range.Text = 'The text discusses the rapid adoption of ChatGPT AI, its limitations compared to human intelligence, ethical concerns, and the impact on the nature of work.' + [char]11 + [char]11 + 'Widely covered aspects include:' + [char]11 + '1. Rapid growth in user base and adoption rate of ChatGPT.' + [char]11 + '2. Predictability of human behavior due to increased internet and smartphone usage, leading to improved AI predictions.' + [char]11 + '3. Limitations of AI tools like ChatGPT in terms of reasoning and creativity.' + [char]11 + '4. Ethical concerns such as reinforcement of ideologies, worldviews, truths, and untruths.' + [char]11 + '5. The potential for AI tools to replace or outsource human work but not change the nature of work fundamentally.' + [char]11 + '6. The need for caution in incorporating AI into working models due to the aggressive push by the technology sector.' + [char]11 + '7. The importance of collaboration between humans and AI rather than a complete takeover or outsourcing.' + [char]11 + '8. The need for society to reflect on the changing nature of human intelligence in the technology age.' + [char]11 + [char]11 + 'Aspects that are being ignored (not mentioned in the text) include:' + [char]11 + '1. Specific examples of how ChatGPT is being used in different industries or sectors.' + [char]11 + '2. Potential benefits of AI beyond efficiency, such as automation of mundane tasks freeing up human labor for more creative work.' + [char]11 + '3. Strategies to ensure ethical development and use of AI, especially in light of recent developments like the dismissal of Microsoft''s AI ethics team.' + [char]11 + '4. Technological advancements or innovations in the field of AI beyond ChatGPT.' + [char]11 + '5. The potential impact of AI on global politics, social dynamics, and cultural norms.'

# Question 3 (paragraph 11)
$d.Paragraphs(11).Range.Text = 'Not mentioned'

# Question 4 (paragraph 13)
$d.Paragraphs(13).Range.Text = 'The article emphasizes that AI tools, like ChatGPT, while powerful and efficient, still lack the ability to reason like humans. The author warns against a complete embrace of these technologies due to ethical concerns, such as their potential to reinforce ideologies and lack of balance between creativity and constraint. The message is to collaborate with AI tools instead of outsourcing or fully replacing human intelligence, and to reflect on how technology is impacting society.'

# Sentiment (paragraph 15)
$d.Paragraphs(15).Range.Text = 'The sentiment is 3.0'

# Entities (paragraph 17)
$d.Paragraphs(17).Range.Text = 'Copyright Syndication Bureau, TikTok, OpenAI, iPhone, Chomsky, AI, ibnezra lang, JOSEPH DANA, Gmail, Joseph Dana, The New York Times Whereas, stu, Twitter ibnezra https twitter.com, ChatGPT-4 https cdn.openai.com, Googles, Arab News, Exponential View, Noam Chomsky'

Write-Host "Done applying edits"
